$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-160 down to 52-161.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly report entry.
# Categorical columns mirror the (now shifted) row 52 below, only the
# measurement columns (date, volume, min/max/avg price, price per kg) change.
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = "2021-12-06"
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112006
$ws.Range("G51").Value = "Repollo"
$ws.Range("H51").Value = "Crespo record"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 200
$ws.Range("K51").Value = 600
$ws.Range("L51").Value = 700
$ws.Range("M51").Value = 650
$ws.Range("N51").Value = "$/unidad"
$ws.Range("O51").Value = "Provincia de Diguillín"
$ws.Range("P51").Value = 650
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"
